$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new Price (column D) value, written as text to preserve
# the original formatting (leading/trailing zeros) of the source data.
$priceUpdates = @{
    2 = "265.08"
    4 = "6.225"
    5 = "0.06156"
    6 = "3.585"
    7 = "6.694"
    8 = "1.362"
    9 = "0.8281"
    10 = "0.01357"
    11 = "0.1615"
    12 = "0.08175"
    13 = "0.03393"
    14 = "0.03145"
    15 = "0.09250"
    16 = "3.908"
    17 = "0.001723"
    18 = "0.04805"
    19 = "0.006287"
    20 = "0.005918"
    21 = "0.001103"
    22 = "0.0001501"
    23 = "3.763"
    24 = "2.301"
    25 = "0.3339"
    26 = "0.1237"
    27 = "0.0002683"
    40 = "0.04620"
    41 = "0.006978"
    42 = "0.1135"
    43 = "0.003402"
    44 = "0.01044"
    45 = "0.00006162"
    47 = "0.7705"
    48 = "0.2045"
    49 = "0.00001401"
    50 = "0.01241"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.Value = "'" + $priceUpdates[$row]
    $cell.Style = "Normal"
}

Write-Output "Updated $($priceUpdates.Count) price cells"
